$d = $word.ActiveDocument

# --- Step 1: remove the "Meta description" paragraph that follows the title ---
$metaPara = $d.Paragraphs.Item(2)
if ($metaPara.Range.Text -like "Meta description*") {
    [void]$metaPara.Range.Delete()
}

# --- Step 2: insert a new bold "Play Big Bucks Bandits Megaways Free | Review 2021"
#             paragraph right before the final (image-prompt) paragraph ---
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
[void]$lastPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$newRange = $newPara.Range
$newXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Big Bucks Bandits Megaways Free | Review 2021</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
[void]$newRange.InsertXML($newXml)

# --- Step 3: replace the old image-prompt text in the final paragraph with the
#             meta-description copy ---
[void]$d.Content.Find.Execute(
    "Create an exciting feature image for Big Bucks Bandits Megaways that will catch the attention of online slot players. The image should be in a cartoon style featuring a happy Maya warrior with glasses. The warrior should be dressed in traditional clothing, with a headdress adorned with gold. In the background, there should be a desert landscape with the Grand Canyon visible. The warrior should be holding a winning slot combination of symbols (such as a horse, guns, cacti, and the like) with an excited expression on their face. The image should be vibrant and colorful to capture the excitement and spirit of the Wild West. This image will entice players to try their luck at Big Bucks Bandits Megaways and experience the thrill of hitting the jackpot.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our unbiased review of Big Bucks Bandits Megaways slot game with 117,649 ways to win. Get ready to play for free!",
    2
)
